$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Split the closing paragraph right after "...though." so that the
#    trailing bookmark (_GoBack) ends up alone in its own paragraph,
#    exactly like the target XML shows (the run text itself is left
#    untouched - only a new paragraph mark is introduced after it).
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute( `
    "redundant in setting the baud rate though.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.Text = "`r"

# ------------------------------------------------------------------
# 2. Insert a blank paragraph followed by the "Command IDs" Heading 1
#    paragraph, right before the (now isolated) bookmark paragraph.
# ------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Last
$insPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$insPoint.InsertBefore("`rCommand IDs`r")

$headingPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$headingPara.Style = "Heading 1"

# ------------------------------------------------------------------
# 3. Insert the "Command IDs" table right before the bookmark
#    paragraph (i.e. directly after the heading).
# ------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Last
$tbl = $d.Tables.Add($bookmarkPara.Range, 4, 2)
$tbl.Style = "Table Grid"

$tbl.ApplyStyleHeadingRows = $true
$tbl.ApplyStyleLastRow = $false
$tbl.ApplyStyleFirstColumn = $true
$tbl.ApplyStyleLastColumn = $false
$tbl.ApplyStyleRowBands = $true
$tbl.ApplyStyleColumnBands = $false

# Column widths: 4675 dxa = 233.75 pt
$tbl.Columns(1).Width = 233.75
$tbl.Columns(2).Width = 233.75

$tbl.Cell(1, 1).Range.Text = "Command"
$tbl.Cell(1, 1).Range.Bold = 1
$tbl.Cell(1, 2).Range.Text = "ID"
$tbl.Cell(1, 2).Range.Bold = 1

$tbl.Cell(2, 1).Range.Text = "Request IMU Data"
$tbl.Cell(2, 2).Range.Text = "102"

$tbl.Cell(3, 1).Range.Text = "Orientation Information"
$tbl.Cell(3, 2).Range.Text = "108"

$tbl.Cell(4, 1).Range.Text = "Send roll, pitch, yaw, and throttle"
$tbl.Cell(4, 2).Range.Text = "200"

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("Tables.Count=" + $d.Tables.Count)
